$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update hours for "Investigación y documentación app" (E6) and "Creación app" (E7)
$ws.Range("E6").Value = 19
$ws.Range("E7").Value = 42

# Recalculate so the SUM formula in E15 reflects the new totals
$excel.Calculate()

# Move the active cell selection to E8 (previously E7)
$ws.Range("E8").Select()
